$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.873.83'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '1.938.21'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4897'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2950'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06901'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '105.10'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07789'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.935.45'
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.351'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7021'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.71%  '
$ws.Range("D17").Value = '30.853.61'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007724'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.574'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.184.86'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9981'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.545'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.869'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.06%  '
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.162'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1041'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.392'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.561'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.571'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.381'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04894'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7627'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.152'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02011'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.659'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.66%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.518'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.093'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9060'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4447'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.712'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.06%  '
$ws.Range("D49").Value = '1.004.75'
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1251'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.53%  '
